{"js": "const body = context.document.body;\n\n// --- Step 1: insert a brand-new list paragraph before the current first\n// paragraph, duplicating the original \"Acoustic communications...\" text,\n// but underlined. insertParagraph() inherits the list style / numbering\n// of the paragraph it is inserted relative to. ---\nlet paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst firstPara = paras.items[0];\nfirstPara.insertParagraph(\n  \"Acoustic communications to evade network security policies (slide 2)\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// Re-load the paragraph collection so we get a \"live\" reference to the\n// freshly inserted paragraph (formatting such as font.underline only\n// reliably stamps both the run AND the paragraph mark once the paragraph\n// has been re-fetched after the sync that created it).\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nparas.items[0].font.underline = Word.UnderlineType.single;\nawait context.sync();\n\n// Re-load again so subsequent index-based access is against the\n// up-to-date/persisted paragraph list.\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// --- Step 2: the texts \"shift\" - each existing paragraph's text is\n// replaced by the text that used to belong to the following paragraph. ---\n\n// old paragraph 0 (\"Acoustic communications...\") -> old paragraph 1's text\nparas.items[1]\n  .getRange(Word.RangeLocation.content)\n  .insertText(\n    \"Retrieve Arxiv redacted data and insight (Arxiv-leaks) (slide 10)\",\n    Word.InsertLocation.replace\n  );\nawait context.sync();\n\n// old paragraph 1 (\"Retrieve Arxiv...\") -> old paragraph 2's text\nparas.items[2]\n  .getRange(Word.RangeLocation.content)\n  .insertText(\n    \"Measurements of Interactions among Android Apps (slide 15)\",\n    Word.InsertLocation.replace\n  );\nawait context.sync();\n\n// old paragraph 2 (\"Measurements of Interactions...\") -> old paragraph 3's text\nparas.items[3]\n  .getRange(Word.RangeLocation.content)\n  .insertText(\n    \"Solving the Android Semantic App (slide 15)\",\n    Word.InsertLocation.replace\n  );\nawait context.sync();\n\n// old paragraph 3 (\"Solving the Android Semantic App...\") -> old paragraph 4's\n// text, and this paragraph becomes underlined too.\nparas.items[4]\n  .getRange(Word.RangeLocation.content)\n  .insertText(\n    \"Physical Layer Device Authentication (slide 29)\",\n    Word.InsertLocation.replace\n  );\nawait context.sync();\n\nparas.items[4].font.underline = Word.UnderlineType.single;\nawait context.sync();\n\n// --- Step 3: the old paragraph 4 (\"Physical Layer Device Authentication\"),\n// whose text has now been absorbed by the previous paragraph, is removed.\n// The final paragraph (\"Tracking Users in Wi-fi...\") is left untouched. ---\nparas.items[5].delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Helper: run Find/Replace scoped to a single paragraph's Range so that\n# replacing text that happens to duplicate text elsewhere in the document\n# (which happens here, since content \"rotates\" down through the list)\n# cannot spill over into neighboring paragraphs.\nfunction Replace-InParagraph($paraIndex, $oldText, $newText) {\n    $p = $d.Paragraphs.Item($paraIndex)\n    $r = $p.Range\n    $r.Find.ClearFormatting()\n    $r.Find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll) - scoped range has only 1 match\n    )\n}\n\n# --- Step 1: insert a brand-new list paragraph before the current first\n# paragraph, duplicating the original \"Acoustic communications...\" text,\n# then underline the whole new paragraph (this stamps the underline on\n# both the run and the paragraph mark itself). ---\n$firstPara = $d.Paragraphs.Item(1)\n$firstPara.Range.InsertParagraphBefore()\n$newPara = $d.Paragraphs.Item(1)\n$newPara.Range.Text = \"Acoustic communications to evade network security policies (slide 2)\"\n$newPara.Range.Font.Underline = 1\n\n# After the insertion every pre-existing paragraph shifted down by one:\n#   2 = old #1 \"Acoustic communications...\"\n#   3 = old #2 \"Retrieve Arxiv...\"\n#   4 = old #3 \"Measurements of Interactions...\"\n#   5 = old #4 \"Solving the Android Semantic App...\"\n#   6 = old #5 \"Physical Layer Device Authentication...\"\n#   7 = old #6 \"Tracking Users in Wi-fi...\" (unchanged)\n\n# --- Step 2: each existing paragraph's text is replaced by the text that\n# used to belong to the following paragraph (content \"rotates\" down). ---\nReplace-InParagraph 2 \"Acoustic communications to evade network security policies (slide 2)\" \"Retrieve Arxiv redacted data and insight (Arxiv-leaks) (slide 10)\"\nReplace-InParagraph 3 \"Retrieve Arxiv redacted data and insight (Arxiv-leaks) (slide 10)\" \"Measurements of Interactions among Android Apps (slide 15)\"\nReplace-InParagraph 4 \"Measurements of Interactions among Android Apps (slide 15)\" \"Solving the Android Semantic App (slide 15)\"\nReplace-InParagraph 5 \"Solving the Android Semantic App (slide 15)\" \"Physical Layer Device Authentication (slide 29)\"\n\n# That last paragraph also becomes underlined.\n$d.Paragraphs.Item(5).Range.Font.Underline = 1\n\n# --- Step 3: the old paragraph (\"Physical Layer Device Authentication\"),\n# whose text has now been absorbed by the previous paragraph, is removed.\n# The final paragraph (\"Tracking Users in Wi-fi...\") is left untouched. ---\n$d.Paragraphs.Item(6).Range.Delete()\n"}
